$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 423.42856
$ws.Range("I4").Value = 93
$ws.Range("J4").Value = 864
$ws.Range("K4").Value = 93
$ws.Range("L4").Value = 864
$ws.Range("M4").Value = 21
$ws.Range("N4").Value = -1092
$ws.Range("H28").Value = 653.5925999999999
$ws.Range("I28").Value = 616.8261
$ws.Range("J28").Value = 865
$ws.Range("K28").Value = 616.8261
$ws.Range("L28").Value = 865
$ws.Range("M28").Value = -131.8261
$ws.Range("N28").Value = -1835
$ws.Range("H51").Value = 16829.312
$ws.Range("I51").Value = 17032.777
$ws.Range("K51").Value = 17032.777
$ws.Range("M51").Value = -16548.777
$ws.Range("H53").Value = 1296.6666
$ws.Range("I53").Value = 300
$ws.Range("J53").Value = 1795
$ws.Range("K53").Value = 300
$ws.Range("L53").Value = 1795
$ws.Range("M53").Value = 337
$ws.Range("N53").Value = -3069
$ws.Range("H62").Value = 7995.9287
$ws.Range("I62").Value = 11421.75
$ws.Range("K62").Value = 11421.75
$ws.Range("M62").Value = -10797.75
$ws.Range("H65").Value = 7995.9287
$ws.Range("I65").Value = 11421.75
$ws.Range("K65").Value = 57108.75
$ws.Range("M65").Value = -53988.75
$ws.Range("H74").Value = 19662.143
$ws.Range("I74").Value = 5911.5
$ws.Range("K74").Value = 5911.5
$ws.Range("M74").Value = -4975.5
$ws.Range("H77").Value = 19662.143
$ws.Range("I77").Value = 5911.5
$ws.Range("K77").Value = 29557.5
$ws.Range("M77").Value = -24877.5
$ws.Range("H101").Value = 11489.444
$ws.Range("I101").Value = 17027
$ws.Range("J101").Value = 414.33334
$ws.Range("K101").Value = 51081
$ws.Range("L101").Value = 1243.00002
$ws.Range("M101").Value = -49459
$ws.Range("N101").Value = -4487.000019999999
$ws.Range("H125").Value = 2684.2144
$ws.Range("I125").Value = 2759.2
$ws.Range("J125").Value = 2642.5557
$ws.Range("K125").Value = 24832.8
$ws.Range("L125").Value = 23783.0013
$ws.Range("M125").Value = -22372.8
$ws.Range("N125").Value = -28703.0013
$ws.Range("H129").Value = 1127.9412
$ws.Range("I129").Value = 840.625
$ws.Range("K129").Value = 2521.875
$ws.Range("M129").Value = 2478.125
$ws.Range("H132").Value = 5570.1113
$ws.Range("I132").Value = 4687.033
$ws.Range("K132").Value = 14061.099
$ws.Range("M132").Value = -11531.099
$ws.Range("H137").Value = 1423.2916
$ws.Range("I137").Value = 1006.6923
$ws.Range("J137").Value = 1915.6364
$ws.Range("K137").Value = 3020.0769
$ws.Range("L137").Value = 5746.9092
$ws.Range("M137").Value = -470.0769
$ws.Range("N137").Value = -10846.9092
$ws.Range("H138").Value = 2276.923
$ws.Range("I138").Value = 1746.25
$ws.Range("J138").Value = 2731.7856
$ws.Range("K138").Value = 5238.75
$ws.Range("L138").Value = 8195.356800000001
$ws.Range("M138").Value = -98.75
$ws.Range("N138").Value = -18475.3568
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 3151.8333
$ws.Range("I110").Value = 2782.2
$ws.Range("K110").Value = 2782.2
$ws.Range("M110").Value = -737.1999999999998
$ws.Range("H132").Value = 29259.902
$ws.Range("I132").Value = 65233.625
$ws.Range("K132").Value = 195700.875
$ws.Range("M132").Value = -193170.875
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 143430.72
$ws.Range("I22").Value = 143430.72
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 143430.72
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -143257.72
$ws.Range("N22").ClearContents()
$ws.Range("H94").Value = 3642.7144
$ws.Range("I94").Value = 3375
$ws.Range("J94").Value = 3999.6667
$ws.Range("K94").Value = 3375
$ws.Range("L94").Value = 3999.6667
$ws.Range("M94").Value = -2924
$ws.Range("N94").Value = -4901.6667
$ws.Range("H105").Value = 2734.5
$ws.Range("I105").Value = 2714.077
$ws.Range("K105").Value = 2714.077
$ws.Range("M105").Value = -967.0770000000002
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 4664.615
$ws.Range("I22").Value = 890
$ws.Range("J22").Value = 9811.817999999999
$ws.Range("K22").Value = 890
$ws.Range("L22").Value = 9811.817999999999
$ws.Range("M22").Value = -540
$ws.Range("N22").Value = -10511.818
$ws.Range("H88").Value = 15335.5
$ws.Range("J88").Value = 15335.5
$ws.Range("L88").Value = 15335.5
$ws.Range("N88").Value = -16147.5
$ws.Range("H91").Value = 15335.5
$ws.Range("J91").Value = 15335.5
$ws.Range("L91").Value = 15335.5
$ws.Range("N91").Value = -18143.5
$ws.Range("H92").Value = 42998
$ws.Range("J92").Value = 42998
$ws.Range("L92").Value = 42998
$ws.Range("N92").Value = -47990
$ws.Range("H94").Value = 2580
$ws.Range("J94").Value = 2546.1667
$ws.Range("L94").Value = 2546.1667
$ws.Range("N94").Value = -3448.1667
$ws.Range("H122").Value = 1684.1538
$ws.Range("I122").Value = 1649.4
$ws.Range("K122").Value = 4948.200000000001
$ws.Range("M122").Value = -2498.200000000001
$ws.Range("H134").Value = 45749.957
$ws.Range("I134").Value = 47478.22
$ws.Range("J134").Value = 6000
$ws.Range("K134").Value = 142434.66
$ws.Range("L134").Value = 18000
$ws.Range("M134").Value = -139899.66
$ws.Range("N134").Value = -23070
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I5").Value = 442.15384
$ws.Range("J5").Value = 1008.2727
$ws.Range("K5").Value = 1326.46152
$ws.Range("L5").Value = 3024.8181
$ws.Range("M5").Value = -1214.46152
$ws.Range("N5").Value = -3248.8181
$ws.Range("H13").Value = 50024.5
$ws.Range("I13").Value = 50024.5
$ws.Range("K13").Value = 150073.5
$ws.Range("M13").Value = -149905.5
$ws.Range("H17").Value = 59.76923
$ws.Range("J17").Value = 99.333336
$ws.Range("L17").Value = 298.000008
$ws.Range("N17").Value = -636.000008
$ws.Range("H113").Value = 1727
$ws.Range("J113").Value = 4888.1113
$ws.Range("L113").Value = 14664.3339
$ws.Range("N113").Value = -19004.3339
$ws.Range("I135").Value = 442.15384
$ws.Range("J135").Value = 1008.2727
$ws.Range("K135").Value = 3979.38456
$ws.Range("L135").Value = 9074.454299999999
$ws.Range("M135").Value = -1444.38456
$ws.Range("N135").Value = -14144.4543
$ws.Range("H140").Value = 1986.0385
$ws.Range("I140").Value = 1608.5
$ws.Range("K140").Value = 4825.5
$ws.Range("M140").Value = 354.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 9058.714
$ws.Range("I13").Value = 2135
$ws.Range("J13").Value = 14251.5
$ws.Range("K13").Value = 2135
$ws.Range("L13").Value = 14251.5
$ws.Range("M13").Value = -1996
$ws.Range("N13").Value = -14529.5
$ws.Range("H122").Value = 3169.5217
$ws.Range("I122").Value = 2176.5881
$ws.Range("J122").Value = 5982.8335
$ws.Range("K122").Value = 6529.7643
$ws.Range("L122").Value = 17948.5005
$ws.Range("M122").Value = -4079.7643
$ws.Range("N122").Value = -22848.5005
$ws.Range("H132").Value = 29843.25
$ws.Range("I132").Value = 39195.965
$ws.Range("K132").Value = 117587.895
$ws.Range("M132").Value = -115057.895
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6676.6924
$ws.Range("I7").Value = 6130.8887
$ws.Range("J7").Value = 7904.75
$ws.Range("K7").Value = 6130.8887
$ws.Range("L7").Value = 7904.75
$ws.Range("M7").Value = -6018.8887
$ws.Range("N7").Value = -8128.75
$ws.Range("H16").Value = 3236.9395
$ws.Range("I16").Value = 2630.6667
$ws.Range("J16").Value = 5965.1665
$ws.Range("K16").Value = 2630.6667
$ws.Range("L16").Value = 5965.1665
$ws.Range("M16").Value = -2460.6667
$ws.Range("N16").Value = -6305.1665
$ws.Range("H22").Value = 88234.42999999999
$ws.Range("J22").Value = 1454
$ws.Range("L22").Value = 1454
$ws.Range("N22").Value = -2044
$ws.Range("H27").Value = 88234.42999999999
$ws.Range("J27").Value = 1454
$ws.Range("L27").Value = 1454
$ws.Range("N27").Value = -1668
$ws.Range("H126").Value = 6676.6924
$ws.Range("I126").Value = 6130.8887
$ws.Range("J126").Value = 7904.75
$ws.Range("K126").Value = 18392.6661
$ws.Range("L126").Value = 23714.25
$ws.Range("M126").Value = -15922.6661
$ws.Range("N126").Value = -28654.25
$ws.Range("H132").Value = 34978
$ws.Range("I132").Value = 46248.223
$ws.Range("K132").Value = 138744.669
$ws.Range("M132").Value = -136214.669
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1407
$ws.Range("I81").Value = 1407
$ws.Range("K81").Value = 2814
$ws.Range("M81").Value = -1753
$ws.Range("H84").Value = 1407
$ws.Range("I84").Value = 1407
$ws.Range("K84").Value = 14070
$ws.Range("M84").Value = -8766
$ws.Range("H122").Value = 6046.1665
$ws.Range("I122").Value = 5555.5
$ws.Range("K122").Value = 16666.5
$ws.Range("M122").Value = -14216.5
$ws.Range("H126").Value = 206060.8
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
